# Updated cryptos list (price "D" column and 1h volume "E" column) per
# the latest GitHub Actions scrape.
#
# Note: several "Price" strings look numeric (e.g. "211.64", "1.00").
# The source workbook stores these as plain text (inline strings), so a
# leading apostrophe ('') is used below for those values to keep Excel
# from auto-converting them to numbers (which would also drop trailing
# zeros, e.g. "1.00" -> 1). Values with two dots (e.g. "27.917.22") are
# already safe since Excel can't parse them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.917.22'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.635.15'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''211.64'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''23.39'
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = '1.867.21'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").Value = '1.643.19'
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = '''65.30'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '27.935.81'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '''229.80'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("E19").Value = '  +3.79%  '
$ws.Range("E20").Value = '  -0.39%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = '''4.36'
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("D23").Value = '''10.14'
$ws.Range("E23").Value = '  -2.31%  '
$ws.Range("D24").Value = '''2.07'
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = '''156.03'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").Value = '''6.98'
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").Value = '''15.56'
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").Value = '''3.10'
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("D34").Value = '1.402.58'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("E35").Value = '  +2.79%  '
$ws.Range("D36").Value = '''1.02'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("E38").Value = '  +0.56%  '
$ws.Range("D39").Value = '''0.559'
$ws.Range("E39").Value = '  -0.49%  '
$ws.Range("D40").Value = '''0.852'
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("D43").Value = '''1.84'
$ws.Range("E43").Value = '  +1.24%  '
$ws.Range("D44").Value = '''66.09'
$ws.Range("E44").Value = '  -1.35%  '
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("D46").Value = '1.775.88'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("D48").Value = '''88.52'
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("D49").Value = '''0.102'
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("D51").Value = '''7.62'
$ws.Range("E51").Value = '  +1.83%  '
